$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Faculty"

# New data row (row 6) - faculty entry (order matters for shared-string allocation)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Sayan Basak"
$ws.Range("E6").Value = "AI"
$ws.Range("F6").Value = "Active"

# New header cells (row 1)
$ws.Range("G1").Value = "Bank Name"
$ws.Range("H1").Value = "Bank A/C No"
$ws.Range("I1").Value = "IFSC Code"

$ws.Range("G6").Value = "Bank Of India"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "402910110001569"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "BKID0004029"

# Bold + text format for the new H1/I1 headers
$ws.Range("H1").NumberFormat = "@"
$ws.Range("I1").NumberFormat = "@"
$ws.Range("H1").Font.Bold = $true
$ws.Range("I1").Font.Bold = $true

# Column widths for new columns (closest reproducible values to the target
# bestFit widths of 11.88671875 / 16.109375 / 11.88671875 given this engine's
# column-width quantization)
$ws.Columns("G").ColumnWidth = 11
$ws.Columns("H").ColumnWidth = 15.333333333333334
$ws.Columns("I").ColumnWidth = 11

# Selection matches the diff
$ws.Range("I7").Select() | Out-Null
